{"js": "// Update each three-digit-by-one-digit division \"problem = quotient, remainder\"\n// line in the table with its new value (commit \"Update master to output\n// generated at c8c62b6\"). Every line of text in the table is unique, so an\n// exact, case-sensitive text search safely targets the right run.\nconst replacements = [\n  [\"598\u00f74=149, 2\", \"934\u00f79=103, 7\"],\n  [\"561\u00f76=93, 3\", \"491\u00f72=245, 1\"],\n  [\"505\u00f78=63, 1\", \"257\u00f72=128, 1\"],\n  [\"587\u00f79=65, 2\", \"380\u00f76=63, 2\"],\n  [\"134\u00f73=44, 2\", \"751\u00f79=83, 4\"],\n  [\"267\u00f74=66, 3\", \"623\u00f78=77, 7\"],\n  [\"291\u00f76=48, 3\", \"515\u00f72=257, 1\"],\n  [\"932\u00f76=155, 2\", \"581\u00f74=145, 1\"],\n  [\"239\u00f79=26, 5\", \"299\u00f76=49, 5\"],\n  [\"900\u00f72=450, 0\", \"986\u00f77=140, 6\"],\n  [\"917\u00f74=229, 1\", \"941\u00f74=235, 1\"],\n  [\"794\u00f78=99, 2\", \"948\u00f79=105, 3\"],\n  [\"328\u00f79=36, 4\", \"122\u00f77=17, 3\"],\n  [\"604\u00f74=151, 0\", \"490\u00f73=163, 1\"],\n  [\"843\u00f77=120, 3\", \"816\u00f74=204, 0\"],\n  [\"589\u00f76=98, 1\", \"879\u00f73=293, 0\"],\n  [\"607\u00f78=75, 7\", \"249\u00f74=62, 1\"],\n  [\"919\u00f73=306, 1\", \"897\u00f77=128, 1\"],\n  [\"729\u00f72=364, 1\", \"477\u00f72=238, 1\"],\n  [\"727\u00f79=80, 7\", \"683\u00f77=97, 4\"],\n  [\"488\u00f78=61, 0\", \"354\u00f73=118, 0\"],\n  [\"918\u00f72=459, 0\", \"224\u00f75=44, 4\"],\n  [\"926\u00f73=308, 2\", \"721\u00f72=360, 1\"],\n  [\"577\u00f73=192, 1\", \"717\u00f79=79, 6\"],\n  [\"274\u00f77=39, 1\", \"541\u00f75=108, 1\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`Text not found: ${oldText}`);\n  }\n\n  results.items[0].insertText(newText, \"Replace\");\n}\n\nawait context.sync();\n", "ps1": "# Update each three-digit-by-one-digit division \"problem = quotient, remainder\"\n# line in the table with its new value (commit \"Update master to output\n# generated at c8c62b6\"). Every line of text in the table is unique, so an\n# exact, case-sensitive Find/Replace safely targets the right run and keeps\n# its original run formatting (font/size) untouched.\n\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"598\u00f74=149, 2\", \"934\u00f79=103, 7\"),\n    @(\"561\u00f76=93, 3\", \"491\u00f72=245, 1\"),\n    @(\"505\u00f78=63, 1\", \"257\u00f72=128, 1\"),\n    @(\"587\u00f79=65, 2\", \"380\u00f76=63, 2\"),\n    @(\"134\u00f73=44, 2\", \"751\u00f79=83, 4\"),\n    @(\"267\u00f74=66, 3\", \"623\u00f78=77, 7\"),\n    @(\"291\u00f76=48, 3\", \"515\u00f72=257, 1\"),\n    @(\"932\u00f76=155, 2\", \"581\u00f74=145, 1\"),\n    @(\"239\u00f79=26, 5\", \"299\u00f76=49, 5\"),\n    @(\"900\u00f72=450, 0\", \"986\u00f77=140, 6\"),\n    @(\"917\u00f74=229, 1\", \"941\u00f74=235, 1\"),\n    @(\"794\u00f78=99, 2\", \"948\u00f79=105, 3\"),\n    @(\"328\u00f79=36, 4\", \"122\u00f77=17, 3\"),\n    @(\"604\u00f74=151, 0\", \"490\u00f73=163, 1\"),\n    @(\"843\u00f77=120, 3\", \"816\u00f74=204, 0\"),\n    @(\"589\u00f76=98, 1\", \"879\u00f73=293, 0\"),\n    @(\"607\u00f78=75, 7\", \"249\u00f74=62, 1\"),\n    @(\"919\u00f73=306, 1\", \"897\u00f77=128, 1\"),\n    @(\"729\u00f72=364, 1\", \"477\u00f72=238, 1\"),\n    @(\"727\u00f79=80, 7\", \"683\u00f77=97, 4\"),\n    @(\"488\u00f78=61, 0\", \"354\u00f73=118, 0\"),\n    @(\"918\u00f72=459, 0\", \"224\u00f75=44, 4\"),\n    @(\"926\u00f73=308, 2\", \"721\u00f72=360, 1\"),\n    @(\"577\u00f73=192, 1\", \"717\u00f79=79, 6\"),\n    @(\"274\u00f77=39, 1\", \"541\u00f75=108, 1\"),\n)\n\n$wdReplaceAll = 2\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.ClearFormatting()\n    $find.Replacement.Text = $newText\n    $find.Forward = $true\n    $find.Wrap = 0\n    $find.Format = $false\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n\n    $found = $find.Execute($oldText, $false, $false, $false, $false, $false, $true, 0, $false, $newText, $wdReplaceAll)\n    if (-not $found) {\n        throw \"Text not found: $oldText\"\n    }\n}\n"}
